$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("initial")
$ws2 = $wb.Worksheets.Item("line_imp")

# Added reactive power limits: new q_lim column on the "initial" sheet
$ws1.Range("G1").Value = "q_lim"
$ws1.Range("G2").Value = 0.5

# Update formula on "line_imp" sheet E3 (NR_iterate_loop_qlim halves the time step again)
$ws2.Range("E3").Formula = "=0.05/2/2"

# Update selection on the "line_imp" sheet
$ws2.Range("E6").Select()

# Restore "initial" as the active sheet/selection
$ws1.Select()
$ws1.Range("G3").Select()
